$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (shifts old D "Tipo" column to E)
$ws.Range("D1").EntireColumn.Insert()

# Fill in the new MAE column header and value
$ws.Range("D1").Value = "MAE"
$ws.Range("D2").Value = 0.8956978851546247
